$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyDesignEstimands")
$ws.Activate()
